# Insert a new weekly data row into the Albahaca price sheet.
# A new row of data (for the most recent week) is inserted just above the
# existing row 110, pushing the old rows 110-142 down to 111-143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 110 (shifts rows 110:142 -> 111:143).
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new week's data.
$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = 44841
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112052
$ws.Range("G110").Value = "Albahaca"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 80
$ws.Range("K110").Value = 8000
$ws.Range("L110").Value = 8000
$ws.Range("M110").Value = 8000
$ws.Range("N110").Value = "`$/paquete"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 8000
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = "Hortaliza"
